# Commit: "Refined metadata to be additional tab"
#
# The "data" sheet's per-row query timestamps (column F) are refreshed to a
# later capture run, and a new "metadata" tab is appended that captures the
# panelapp query metadata (panel name/id/version/version-created timestamp/
# query time/request URL) that used to live elsewhere.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. Refresh the "time_taken" timestamps on the existing "data" sheet
# ---------------------------------------------------------------------
$ws1.Range("F2").Value  = "2021-10-05 14:35:17.706406"
$ws1.Range("F3").Value  = "2021-10-05 14:35:17.706414"
$ws1.Range("F4").Value  = "2021-10-05 14:35:17.706417"
$ws1.Range("F5").Value  = "2021-10-05 14:35:17.706419"
$ws1.Range("F6").Value  = "2021-10-05 14:35:17.706422"
$ws1.Range("F7").Value  = "2021-10-05 14:35:17.706425"
$ws1.Range("F8").Value  = "2021-10-05 14:35:17.706427"
$ws1.Range("F9").Value  = "2021-10-05 14:35:17.706430"
$ws1.Range("F10").Value = "2021-10-05 14:35:17.706433"
$ws1.Range("F11").Value = "2021-10-05 14:35:17.706435"
$ws1.Range("F12").Value = "2021-10-05 14:35:17.706438"
$ws1.Range("F13").Value = "2021-10-05 14:35:17.706440"
$ws1.Range("F14").Value = "2021-10-05 14:35:17.706443"
$ws1.Range("F15").Value = "2021-10-05 14:35:17.706445"
$ws1.Range("F16").Value = "2021-10-05 14:35:17.706447"
$ws1.Range("F17").Value = "2021-10-05 14:35:17.706450"
$ws1.Range("F18").Value = "2021-10-05 14:35:17.706453"
$ws1.Range("F19").Value = "2021-10-05 14:35:17.706455"
$ws1.Range("F20").Value = "2021-10-05 14:35:17.706457"
$ws1.Range("F21").Value = "2021-10-05 14:35:17.706460"
$ws1.Range("F22").Value = "2021-10-05 14:35:17.706462"
$ws1.Range("F23").Value = "2021-10-05 14:35:17.706465"
$ws1.Range("F24").Value = "2021-10-05 14:35:17.706467"
$ws1.Range("F25").Value = "2021-10-05 14:35:17.706470"
$ws1.Range("F26").Value = "2021-10-05 14:35:17.706472"
$ws1.Range("F27").Value = "2021-10-05 14:35:17.706475"
$ws1.Range("F28").Value = "2021-10-05 14:35:17.706477"
$ws1.Range("F29").Value = "2021-10-05 14:35:17.706480"
$ws1.Range("F30").Value = "2021-10-05 14:35:17.706482"
$ws1.Range("F31").Value = "2021-10-05 14:35:17.706485"
$ws1.Range("F32").Value = "2021-10-05 14:35:17.706487"
$ws1.Range("F33").Value = "2021-10-05 14:35:17.706490"
$ws1.Range("F34").Value = "2021-10-05 14:35:17.706492"
$ws1.Range("F35").Value = "2021-10-05 14:35:17.706495"
$ws1.Range("F36").Value = "2021-10-05 14:35:17.706497"
$ws1.Range("F37").Value = "2021-10-05 14:35:17.706500"
$ws1.Range("F38").Value = "2021-10-05 14:35:17.706502"
$ws1.Range("F39").Value = "2021-10-05 14:35:17.706505"
$ws1.Range("F40").Value = "2021-10-05 14:35:17.706507"
$ws1.Range("F41").Value = "2021-10-05 14:35:17.706510"
$ws1.Range("F42").Value = "2021-10-05 14:35:17.706513"
$ws1.Range("F43").Value = "2021-10-05 14:35:17.706515"
$ws1.Range("F44").Value = "2021-10-05 14:35:17.706518"
$ws1.Range("F45").Value = "2021-10-05 14:35:17.706520"
$ws1.Range("F46").Value = "2021-10-05 14:35:17.706522"
$ws1.Range("F47").Value = "2021-10-05 14:35:17.706525"

# ---------------------------------------------------------------------
# 2. Append a new "metadata" worksheet after the existing "data" sheet
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$meta = $wb.Worksheets.Add($null, $lastSheet)
$meta.Name = "metadata"

# Match the outline summary placement used on the "data" sheet
# (<outlinePr summaryBelow="1" summaryRight="1"/>).
$meta.Outline.SummaryRow = 1
$meta.Outline.SummaryColumn = 1

# Header row (B1:G1)
$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

# Copy the bold/bordered/centred header style already used by sheet "data"
# (row 1 header cells) onto the new header row, so both tabs share the same
# cellXfs entry instead of minting a duplicate style.
$ws1.Range("B1").Copy()
$meta.Range("B1:G1").PasteSpecial(-4122)   # xlPasteFormats

# Data row (A2:G2) -- pandas-style index column in A, data in B..G
$meta.Range("A2").Value = 0
$ws1.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)      # xlPasteFormats (matches "data" index style)

$meta.Range("B2").Value = "Predominantly Antibody Deficiency"
$meta.Range("C2").Value = 222
$meta.Range("E2").Value = "2021-08-17T08:22:45.084260Z"
$meta.Range("F2").Value = "2021-10-05 14:35:17.702538"
$meta.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/222/?format=json"

# D2 ("0.92") must stay a plain text value, not get auto-coerced into the
# number 0.92 -- build it as a text formula in a scratch cell, then paste
# only the *value* (not the format) into D2 so it keeps its string type
# without acquiring a new/duplicate style.
$meta.Range("Z1").Formula = '="0.92"'
$meta.Range("Z1").Copy()
$meta.Range("D2").PasteSpecial(-4163)      # xlPasteValues
$meta.Range("Z1").ClearContents()

# Restore "data" as the active tab (adding a sheet makes the new one
# active by default) so the workbook's selection state is left as it was.
$ws1.Activate()
